# Update "想去人数" (interested count) values in column F across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.
# 本地生活 (Local life) sheet is left untouched.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 81
$ws.Range("F5").Value = 1667
$ws.Range("F6").Value = 3267
$ws.Range("F7").Value = 836
$ws.Range("F8").Value = 2065
$ws.Range("F9").Value = 1982
$ws.Range("F10").Value = 1022
$ws.Range("F11").Value = 356
$ws.Range("F13").Value = 1615
$ws.Range("F18").Value = 80
$ws.Range("F19").Value = 1451
$ws.Range("F20").Value = 532
$ws.Range("F21").Value = 638
$ws.Range("F22").Value = 324
$ws.Range("F23").Value = 10775
$ws.Range("F24").Value = 11754
$ws.Range("F26").Value = 665
$ws.Range("F27").Value = 1842
$ws.Range("F28").Value = 150
$ws.Range("F29").Value = 453

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 34

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 81
$ws.Range("F7").Value = 1667
$ws.Range("F8").Value = 3267
$ws.Range("F9").Value = 836
$ws.Range("F10").Value = 2065
$ws.Range("F11").Value = 1982
$ws.Range("F12").Value = 1022
$ws.Range("F13").Value = 356
$ws.Range("F15").Value = 1615
$ws.Range("F22").Value = 81
$ws.Range("F23").Value = 1451
$ws.Range("F24").Value = 532
$ws.Range("F25").Value = 638
$ws.Range("F26").Value = 324
$ws.Range("F27").Value = 10775
$ws.Range("F28").Value = 11754
$ws.Range("F30").Value = 665
$ws.Range("F31").Value = 1842
$ws.Range("F33").Value = 34
$ws.Range("F34").Value = 150
$ws.Range("F35").Value = 453

$wb.Save()
